$d = $word.ActiveDocument

# Mapping of old text -> new text for this day's worksheet update.
$replacements = @(
    @{old = "2024-05-21 Tuesday"; new = "2024-05-22 Wednesday"},
    @{old = "72×70=5040"; new = "81×58=4698"},
    @{old = "93×94=8742"; new = "45×71=3195"},
    @{old = "77×60=4620"; new = "39×47=1833"},
    @{old = "47×92=4324"; new = "17×33=561"},
    @{old = "45×83=3735"; new = "35×13=455"},
    @{old = "47×57=2679"; new = "65×89=5785"},
    @{old = "48×62=2976"; new = "44×28=1232"},
    @{old = "26×13=338"; new = "14×41=574"},
    @{old = "72×72=5184"; new = "47×29=1363"},
    @{old = "59×67=3953"; new = "16×26=416"},
    @{old = "58×82=4756"; new = "60×28=1680"},
    @{old = "41×37=1517"; new = "52×27=1404"},
    @{old = "85×24=2040"; new = "74×98=7252"},
    @{old = "40×45=1800"; new = "17×50=850"},
    @{old = "13×83=1079"; new = "52×19=988"},
    @{old = "79×62=4898"; new = "95×94=8930"},
    @{old = "87×51=4437"; new = "92×17=1564"},
    @{old = "89×49=4361"; new = "51×27=1377"},
    @{old = "38×97=3686"; new = "20×48=960"},
    @{old = "14×54=756"; new = "18×86=1548"},
    @{old = "42×11=462"; new = "80×26=2080"},
    @{old = "47×13=611"; new = "63×31=1953"},
    @{old = "35×28=980"; new = "82×86=7052"},
    @{old = "46×62=2852"; new = "12×88=1056"},
    @{old = "55×68=3740"; new = "54×11=594"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
